$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right after "2021-Q4" (i.e. right before
#    the existing "总计" summary sheet).
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $afterSheet)
$q1.Name = "2022-Q1"

# --- Header row (B1:H1) -----------------------------------------------------
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1Hdr = $q1.Range("B1:H1")
$q1Hdr.Borders.LineStyle = 1
$q1Hdr.Font.Bold = $true
$q1Hdr.HorizontalAlignment = -4108
$q1Hdr.VerticalAlignment = -4160

# --- Row 2: 009686 / 华夏磐利一年定期开放混合A ------------------------------
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'009686"
$q1.Range("C2").Value = "华夏磐利一年定期开放混合A"
$q1.Range("D2").Value = "'16.02"
$q1.Range("E2").Value = "'93.69"
$q1.Range("F2").Value = "'4.25"
$q1.Range("G2").Value = "'0.6808"
$q1.Range("H2").Value = 4

# --- Row 3: 009687 / 华夏磐利一年定期开放混合C ------------------------------
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'009687"
$q1.Range("C3").Value = "华夏磐利一年定期开放混合C"
$q1.Range("D3").Value = "'0.51"
$q1.Range("E3").Value = "'93.69"
$q1.Range("F3").Value = "'4.25"
$q1.Range("G3").Value = "'0.0217"
$q1.Range("H3").Value = 4

# Column A (row index) carries the same bold/border/centred look as the header
$q1ColA = $q1.Range("A2:A3")
$q1ColA.Borders.LineStyle = 1
$q1ColA.Font.Bold = $true
$q1ColA.HorizontalAlignment = -4108
$q1ColA.VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: add a new top data row for 2022-Q1 and
#    push the previous 2021-Q4 / 2021-Q3 rows down by one.
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Shift the existing two data rows down (row3<-row2, row4<-row3) first so we
# don't clobber data we still need while re-numbering the index column.
$tot.Range("A4").Value = 2
$tot.Range("B4").Value = "2021-Q3"
$tot.Range("C4").Value = 2
$tot.Range("D4").Value = 0.65

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2021-Q4"
$tot.Range("C3").Value = 4
$tot.Range("D3").Value = 0.33

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 2
$tot.Range("D2").Value = 0.7

# Re-apply the bold/border/centred style used by the index column (A) for
# every data row, matching the look already used for row 1 / header.
$totColA = $tot.Range("A2:A4")
$totColA.Borders.LineStyle = 1
$totColA.Font.Bold = $true
$totColA.HorizontalAlignment = -4108
$totColA.VerticalAlignment = -4160

Write-Output "2022-Q1 sheet added and 总计 sheet updated"
